# Release Aspose.Cells Cloud SDK 23.12 - workbook now ships 3 sheets,
# Sheet1 carries a sample value, and the saved page setup matches A4 portrait.

$wb = $excel.ActiveWorkbook

# --- Sheet1: write the sample value -------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = 1111

# --- Add Sheet2 / Sheet3 right after Sheet1, preserving left-to-right order
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# --- Page setup for every sheet (A4 portrait) ----------------------------
foreach ($s in $wb.Worksheets) {
    $s.PageSetup.PaperSize = 9
    $s.PageSetup.Orientation = 1
}

# --- Restore Sheet1 as the active/selected sheet, with A2 selected -------
$ws1.Select() | Out-Null
$ws1.Range("A2").Select() | Out-Null
